# Generate Report for Handoff
#
# A handoff batch was (re)generated, which refreshes the generation/handoff
# timestamps and marks the affected files as hot-fix priority ("ht") on the
# Overview, zh-cn and de-de sheets. The file "cf740957-..." (row 13) was not
# part of this batch and is left untouched.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 11, 12, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
# --- de-de sheet:    column H = "Latest Handoff Datetime" (same value) ---
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-15 18:20:36"
    $wsDeDe.Range("H$r").Value = "2016-08-15 18:20:36"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime" ---
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-15 18:20:31"
}

# --- zh-cn / de-de sheets: column E = "Priority" -> hot-fix ("ht") ---
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
